# Add 5 new user accounts (D001..D005 / pass001..pass005) to the Users sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (207) down onto the
# five new rows (208-212) so the new cells keep the same style (s="2").
$ws.Range("A207:B207").Copy()
$ws.Range("A208:B212").PasteSpecial(-4122)

# Fill in the values. The order below reproduces the order in which the
# author apparently typed them: row 208 and 209 filled in left-to-right
# (username then password), then the usernames for rows 210-212 were
# entered before circling back to fill in their passwords - this matches
# the shared-string insertion order seen in the target workbook.
$ws.Range("A208").Value() = "D001"
$ws.Range("B208").Value() = "pass001"
$ws.Range("A209").Value() = "D002"
$ws.Range("B209").Value() = "pass002"
$ws.Range("A210").Value() = "D003"
$ws.Range("A211").Value() = "D004"
$ws.Range("A212").Value() = "D005"
$ws.Range("B210").Value() = "pass003"
$ws.Range("B211").Value() = "pass004"
$ws.Range("B212").Value() = "pass005"

# Update the selected cell to match where the author ended up after typing.
$ws.Range("E208").Select()
